# Update the "想去人数" (interest count) figures for several events, and the
# "最低票价" (lowest ticket price) for one event, across the "展览", "演出"
# and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 862
$ws.Range("F3").Value  = 1435
$ws.Range("F4").Value  = 1081
$ws.Range("F5").Value  = 509
$ws.Range("F7").Value  = 657
$ws.Range("F8").Value  = 234
$ws.Range("F10").Value = 74
$ws.Range("F11").Value = 212
$ws.Range("F12").Value = 143
$ws.Range("G12").Value = 55
$ws.Range("F13").Value = 1774
$ws.Range("F14").Value = 424
$ws.Range("F16").Value = 486
$ws.Range("F17").Value = 252
$ws.Range("F21").Value = 655
$ws.Range("F22").Value = 45
$ws.Range("F26").Value = 1525
$ws.Range("F27").Value = 269

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 22
$ws.Range("F7").Value = 7

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 862
$ws.Range("F4").Value  = 1435
$ws.Range("F5").Value  = 1081
$ws.Range("F8").Value  = 509
$ws.Range("F10").Value = 657
$ws.Range("F12").Value = 234
$ws.Range("F14").Value = 74
$ws.Range("F15").Value = 212
$ws.Range("F16").Value = 143
$ws.Range("G16").Value = 55
$ws.Range("F17").Value = 1774
$ws.Range("F19").Value = 424
$ws.Range("F21").Value = 486
$ws.Range("F22").Value = 252
$ws.Range("F24").Value = 22
$ws.Range("F27").Value = 7
$ws.Range("F30").Value = 655
$ws.Range("F35").Value = 45
$ws.Range("F39").Value = 1525
$ws.Range("F40").Value = 269
